$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "R4" reference designator to "R4, R6" now that R6's dedicated
# bleeder-resistor line item is being removed and folded into R4's note.
$ws.Range("F25").Value = "R4, R6"

# The bleeder resistor quantity actually needed is 2, not 1 (E25 recalculates
# automatically since it holds the formula =C25*D25).
$ws.Range("D25").Value = 2
$ws.Rows("25").RowHeight = 13.4

# Remove the now-obsolete "Bleeder resistor for leveling capacitor" (R6) line
# item entirely -- it was a separate part; its function is now covered by R4.
# Deleting the row shifts everything below it up by one and fixes up the
# SUM() formula range automatically.
$ws.Rows("29").Delete()

# Materialize the previously-implicit blank spacer row 10 (between the two
# "first batch" part groups) with its normal default row height.
$ws.Rows("10").RowHeight = 12.1
$ws.Rows("10").UseStandardHeight = $true

# Restore the selection to where the user finished editing.
$ws.Range("B38").Select()
